# Box Plot Updates, Color Updates Main Figures
#
# Reposition the clade-label / count textboxes (tx9..tx18) that live
# inside the group shape on slide 1. The OOXML <a:off> values are in
# EMU (914400 EMU = 1 inch = 72 pt, i.e. 12700 EMU per point); the
# PowerPoint COM object model's Shape.Left/Shape.Top are in points, so
# each target EMU value below is pre-divided by 12700 (with the literal
# nudged by <1 EMU of slack to absorb the host's internal float32
# rounding when it re-serializes points back to EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(2)

# tx9  "Stramenopiles "   (6290859,3710194) -> (6291697,3709492)
$g.GroupItems.Item(7).Left = 495.40924072265625
$g.GroupItems.Item(7).Top = 292.08599853515625

# tx10 " 85 / 110"        (6591990,4062841) -> (6592828,4062139)
$g.GroupItems.Item(8).Left = 519.120361328125
$g.GroupItems.Item(8).Top = 319.8534851074219

# tx11 "Rhizaria "        (5093640,5488928) -> (5094137,5489458)
$g.GroupItems.Item(9).Left = 401.1131591796875
$g.GroupItems.Item(9).Top = 432.2408142089844

# tx12 " 6 / 6"           (5280331,5798382) -> (5280828,5798912)
$g.GroupItems.Item(10).Left = 415.813232421875
$g.GroupItems.Item(10).Top = 456.6072692871094

# tx13 "Metamonada "      (3823067,4708492) -> (3822863,4708525)
$g.GroupItems.Item(11).Left = 301.0128479003906
$g.GroupItems.Item(11).Top = 370.75

# tx14 " 15 / 15"         (4124303,5017946) -> (4124100,5017979)
$g.GroupItems.Item(12).Left = 324.7322998046875
$g.GroupItems.Item(12).Top = 395.1164855957031

# tx15 "Discoba "         (3545136,4126561) -> (3543537,4128100)
$g.GroupItems.Item(13).Left = 279.0186767578125
$g.GroupItems.Item(13).Top = 325.0472717285156

# tx16 " 42 / 42"         (3605267,4436015) -> (3603668,4437554)
$g.GroupItems.Item(14).Left = 283.7533874511719
$g.GroupItems.Item(14).Top = 349.4137268066406

# tx17 "Alveolata "       (4909156,2522144) -> (4323598,3107783)
$g.GroupItems.Item(15).Left = 340.4407958984375
$g.GroupItems.Item(15).Top = 244.70733642578125

# tx18 " 17 / 135"        (4963305,2831492) -> (4377748,3417131)
$g.GroupItems.Item(16).Left = 344.70458984375
$g.GroupItems.Item(16).Top = 269.0654602050781
